# Apply cryptos list price/volume/ranking update
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = '60.222.72'
$ws.Range("E2").Value = '  +5.65%  '

# Row 3
$ws.Range("D3").Value = '3.267.56'
$ws.Range("E3").Value = '  +0.96%  '

# Row 4
$c = $ws.Range("D4")
$c.Value = "'" + '0.996'
$c.Style = "Normal"
$ws.Range("E4").Value = '  -0.43%  '

# Row 5
$c = $ws.Range("D5")
$c.Value = "'" + '406.49'
$c.Style = "Normal"
$ws.Range("E5").Value = '  +2.93%  '

# Row 6
$c = $ws.Range("D6")
$c.Value = "'" + '110.63'
$c.Style = "Normal"
$ws.Range("E6").Value = '  +3.12%  '

# Row 7
$ws.Range("D7").Value = '3.411.45'
$ws.Range("E7").Value = '  +5.51%  '

# Row 8
$c = $ws.Range("D8")
$c.Value = "'" + '0.560'
$c.Style = "Normal"
$ws.Range("E8").Value = '  -2.37%  '

# Row 9
$c = $ws.Range("D9")
$c.Value = "'" + '0.999'
$c.Style = "Normal"
$ws.Range("E9").Value = '  -0.06%  '

# Row 10
$c = $ws.Range("D10")
$c.Value = "'" + '0.609'
$c.Style = "Normal"
$ws.Range("E10").Value = '  -1.44%  '

# Row 11
$c = $ws.Range("D11")
$c.Value = "'" + '0.106'
$c.Style = "Normal"
$ws.Range("E11").Value = '  +10.40%  '

# Row 12
$c = $ws.Range("D12")
$c.Value = "'" + '38.07'
$c.Style = "Normal"
$ws.Range("E12").Value = '  -2.54%  '

# Row 13
$ws.Range("E13").Value = '  +0.10%  '

# Row 14
$ws.Range("D14").Value = '3.719.65'
$ws.Range("E14").Value = '  -0.83%  '

# Row 15
$c = $ws.Range("D15")
$c.Value = "'" + '8.06'
$c.Style = "Normal"
$ws.Range("E15").Value = '  -1.82%  '

# Row 16
$c = $ws.Range("D16")
$c.Value = "'" + '18.73'
$c.Style = "Normal"
$ws.Range("E16").Value = '  -2.10%  '

# Row 17
$ws.Range("D17").Value = '3.275.98'
$ws.Range("E17").Value = '  +1.75%  '

# Row 18
$ws.Range("D18").Value = '59.698.46'
$ws.Range("E18").Value = '  +5.04%  '

# Row 19
$c = $ws.Range("D19")
$c.Value = "'" + '0.983'
$c.Style = "Normal"
$ws.Range("E19").Value = '  -5.00%  '

# Row 20
$c = $ws.Range("D20")
$c.Value = "'" + '10.35'
$c.Style = "Normal"
$ws.Range("E20").Value = '  -4.78%  '

# Row 21
$c = $ws.Range("D21")
$c.Value = "'" + '0.0000110'
$c.Style = "Normal"
$ws.Range("E21").Value = '  +3.28%  '

# Row 22
$c = $ws.Range("D22")
$c.Value = "'" + '3.16'
$c.Style = "Normal"
$ws.Range("E22").Value = '  -5.37%  '

# Row 23
$c = $ws.Range("D23")
$c.Value = "'" + '292.60'
$c.Style = "Normal"
$ws.Range("E23").Value = '  -1.36%  '

# Row 24
$c = $ws.Range("D24")
$c.Value = "'" + '12.12'
$c.Style = "Normal"
$ws.Range("E24").Value = '  -6.54%  '

# Row 25
$c = $ws.Range("D25")
$c.Value = "'" + '72.72'
$c.Style = "Normal"
$ws.Range("E25").Value = '  -1.73%  '

# Row 26
$ws.Range("E26").Value = '  -4.62%  '

# Row 27
$c = $ws.Range("D27")
$c.Value = "'" + '4.46'
$c.Style = "Normal"
$ws.Range("E27").Value = '  +2.38%  '

# Row 28
$c = $ws.Range("D28")
$c.Value = "'" + '28.20'
$c.Style = "Normal"
$ws.Range("E28").Value = '  +1.38%  '

# Row 29
$c = $ws.Range("D29")
$c.Value = "'" + '7.29'
$c.Style = "Normal"
$ws.Range("E29").Value = '  +0.47%  '

# Row 30
$ws.Range("B30").Value = 'Kaspa'
$ws.Range("C30").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$c = $ws.Range("D30")
$c.Value = "'" + '0.168'
$c.Style = "Normal"
$ws.Range("E30").Value = '  -0.57%  '

# Row 31
$ws.Range("B31").Value = 'Filecoin'
$ws.Range("C31").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$c = $ws.Range("D31")
$c.Value = "'" + '7.40'
$c.Style = "Normal"
$ws.Range("E31").Value = '  -3.81%  '

# Row 32
$c = $ws.Range("D32")
$c.Value = "'" + '0.998'
$c.Style = "Normal"
$ws.Range("E32").Value = '  -0.15%  '

# Row 33
$c = $ws.Range("D33")
$c.Value = "'" + '10.99'
$c.Style = "Normal"
$ws.Range("E33").Value = '  -3.72%  '

# Row 34
$c = $ws.Range("D34")
$c.Value = "'" + '0.107'
$c.Style = "Normal"
$ws.Range("E34").Value = '  -1.73%  '

# Row 35
$c = $ws.Range("D35")
$c.Value = "'" + '2.38'
$c.Style = "Normal"
$ws.Range("E35").Value = '  +12.66%  '

# Row 36
$c = $ws.Range("D36")
$c.Value = "'" + '38.96'
$c.Style = "Normal"
$ws.Range("E36").Value = '  +2.54%  '

# Row 37
$c = $ws.Range("D37")
$c.Value = "'" + '52.09'
$c.Style = "Normal"
$ws.Range("E37").Value = '  +0.75%  '

# Row 38
$c = $ws.Range("D38")
$c.Value = "'" + '0.0465'
$c.Style = "Normal"
$ws.Range("E38").Value = '  -3.76%  '

# Row 39
$c = $ws.Range("D39")
$c.Value = "'" + '0.994'
$c.Style = "Normal"
$ws.Range("E39").Value = '  -0.56%  '

# Row 40
$c = $ws.Range("D40")
$c.Value = "'" + '3.00'
$c.Style = "Normal"
$ws.Range("E40").Value = '  +1.56%  '

# Row 41
$ws.Range("B41").Value = 'EnergySwap'
$ws.Range("C41").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$c = $ws.Range("D41")
$c.Value = "'" + '26.46'
$c.Style = "Normal"
$ws.Range("E41").Value = '  +19.29%  '

# Row 42
$ws.Range("B42").Value = 'LidoDAOToken'
$ws.Range("C42").Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$c = $ws.Range("D42")
$c.Value = "'" + '3.26'
$c.Style = "Normal"
$ws.Range("E42").Value = '  -7.76%  '

# Row 43
$ws.Range("B43").Value = 'Monero'
$ws.Range("C43").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$c = $ws.Range("D43")
$c.Value = "'" + '134.72'
$c.Style = "Normal"
$ws.Range("E43").Value = '  +0.15%  '

# Row 44
$ws.Range("B44").Value = 'Stellar'
$ws.Range("C44").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$c = $ws.Range("D44")
$c.Value = "'" + '0.118'
$c.Style = "Normal"
$ws.Range("E44").Value = '  -1.84%  '

# Row 45
$ws.Range("B45").Value = 'ARBITRUM'
$ws.Range("C45").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$c = $ws.Range("D45")
$c.Value = "'" + '1.85'
$c.Style = "Normal"
$ws.Range("E45").Value = '  -2.02%  '

# Row 46
$ws.Range("B46").Value = 'TheGraph'
$ws.Range("C46").Value = 'https://coinranking.com/coin/qhd1biQ7M+thegraph-grt'
$c = $ws.Range("D46")
$c.Value = "'" + '0.273'
$c.Style = "Normal"
$ws.Range("E46").Value = '  -2.91%  '

# Row 47
$ws.Range("B47").Value = 'Celestia'
$ws.Range("C47").Value = 'https://coinranking.com/coin/YQcD0lBl7+celestia-tia'
$c = $ws.Range("D47")
$c.Value = "'" + '16.02'
$c.Style = "Normal"
$ws.Range("E47").Value = '  -6.02%  '

# Row 48
$ws.Range("B48").Value = 'NEARProtocol'
$ws.Range("C48").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$c = $ws.Range("D48")
$c.Value = "'" + '3.70'
$c.Style = "Normal"
$ws.Range("E48").Value = '  -6.72%  '

# Row 49
$ws.Range("B49").Value = 'WEMIXToken'
$ws.Range("C49").Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$c = $ws.Range("D49")
$c.Value = "'" + '2.18'
$c.Style = "Normal"
$ws.Range("E49").Value = '  +3.32%  '

# Row 50
$ws.Range("D50").Value = '3.824.72'
$ws.Range("E50").Value = '  +7.59%  '

# Row 51
$ws.Range("D51").Value = '2.095.89'
$ws.Range("E51").Value = '  -2.83%  '
